# The sheet used to hold a single column (A1:A13) with a numeric header (0)
# in A1 and scattered "so"/"maybe"/"no" answers further down the column.
# It is replaced by a compact two-column table (A1:B4): column A keeps the
# header "Unnamed: 0" plus the three answers "so"/"maybe"/"no" on the rows
# right below it, and a new column B ("yes") is added alongside, left blank
# for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: replace the header value and compact the answers up to A2:A4
$ws.Range("A1").Value = "Unnamed: 0"
$ws.Range("A2").Value = "so"
$ws.Range("A3").Value = "maybe"
$ws.Range("A4").Value = "no"

# Drop the rest of the old column (old data ran all the way to A13)
$ws.Range("A5:A13").Clear()

# --- Column B: new "yes" column, header only, data cells blank
$ws.Range("B1").Value = "yes"

# Give the new header cell (B1) the same look as the existing header (A1):
# bold font, thin border, centered/top aligned.
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats

# Touch B2:B4 so they exist as real (blank) cells under the "yes" column,
# matching the blank cells that accompany every data row in column A.
$ws.Cells.Item(2, 2).Font.Italic = $false
$ws.Cells.Item(3, 2).Font.Italic = $false
$ws.Cells.Item(4, 2).Font.Italic = $false
